$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (shifts robinson/wintri rows down to 8/9)
$ws.Rows("7:7").Insert()

# Populate the new "mercator" entry.
# Leading apostrophe forces text storage (matching the quotePrefix style used
# by the other proj4-string cells in column B, since the value starts with "+").
$ws.Range("A7").Value = "mercator"
$ws.Range("B7").Value = "'+proj=merc +lon_0=0 +k=1 +x_0=0 +y_0=0 +ellps=WGS84 +datum=WGS84 +units=m +no_defs"""

# Re-apply the alphabetical sort over the now-larger range (A2:B9) so the
# worksheet's sortState/sortCondition reflect the new extent.
$ws.Sort.SetRange($ws.Range("A2:B9"))
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A9"))
$ws.Sort.Apply()

# Match the author's final selection in the saved file.
$ws.Range("B6").Select()
